# Actualización 11 de Mayo - Tarde
# Updates the "2o Parcial" and "Final" sheets with refreshed exam statistics.

$wb = $excel.ActiveWorkbook

# ---- "2o Parcial" sheet ----
$ws2 = $wb.Worksheets.Item("2o Parcial")

# Row 4 - Avila Coronado Julieta
$ws2.Range("C4").Value = 0
$ws2.Range("D4").Value = 0
$ws2.Range("E4").Value = 117
$ws2.Range("F4").Value = 71.34
$ws2.Range("G4").Value = 47
$ws2.Range("H4").Value = 28.66
$ws2.Range("I4").Value = 7.5

# Row 15 - Fernández Castro Araceli
$ws2.Range("C15").Value = 0
$ws2.Range("D15").Value = 0
$ws2.Range("E15").Value = 21
$ws2.Range("F15").Value = 56.76
$ws2.Range("G15").Value = 16
$ws2.Range("H15").Value = 43.24
$ws2.Range("I15").Value = 6.4

# Row 23 - Hernández Mendoza Delfina
$ws2.Range("C23").Value = 54
$ws2.Range("D23").Value = 47.79
$ws2.Range("E23").Value = 59
$ws2.Range("F23").Value = 52.21
$ws2.Range("G23").Value = 10
$ws2.Range("H23").Value = 8.85
$ws2.Range("I23").Value = 8.5

# Row 47 - Velasco Sanchez David
$ws2.Range("C47").Value = 84
$ws2.Range("D47").Value = 42.64
$ws2.Range("E47").Value = 109
$ws2.Range("F47").Value = 55.33
$ws2.Range("G47").Value = 88
$ws2.Range("H47").Value = 44.67
$ws2.Range("I47").Value = 7.9

# ---- "Final" sheet ----
$ws3 = $wb.Worksheets.Item("Final")

# Row 4 - Avila Coronado Julieta
$ws3.Range("E4").Value = 117
$ws3.Range("F4").Value = 71.34
$ws3.Range("G4").Value = 47
$ws3.Range("H4").Value = 28.66

# Row 15 - Fernández Castro Araceli
$ws3.Range("E15").Value = 21
$ws3.Range("F15").Value = 56.76
$ws3.Range("G15").Value = 16
$ws3.Range("H15").Value = 43.24

# Row 47 - Velasco Sanchez David
$ws3.Range("E47").Value = 140
$ws3.Range("F47").Value = 71.06999999999999
$ws3.Range("G47").Value = 57
$ws3.Range("H47").Value = 28.93
